$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new task rows (row 4 -> task 3, row 5 -> task 4)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Setup fast api"
$ws.Range("C4").Value = "Create a dummy end point using fast api"
$ws.Range("D4").Value = "Not Started"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Setup swagger page"
$ws.Range("C5").Value = "Setup swagger page for endpoints"
$ws.Range("D5").Value = "Not Started"

# Columns B and C auto-fit wider to accommodate the new, longer text
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 34

# Match the recorded selection state after the edit
$ws.Range("D7").Select()
